# Add "I0" (column I) and "IF" (column J) headers plus their per-row values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1), styled like the other header cells (bold/border/center).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Per-row data for columns I (I0) and J (IF), rows 2-38.
$data = @{
    2  = @(8, 9)
    3  = @(8, 9)
    4  = @(9, 9)
    5  = @(9, 9)
    6  = @(9, 9)
    7  = @(9, 9)
    8  = @(9, 9)
    9  = @(8, 8)
    10 = @(9, 9)
    11 = @(5, 5)
    12 = @(5, 5)
    13 = @(8, 8)
    14 = @(11, 11)
    15 = @(7, 7)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(6, 6)
    19 = @(5, 6)
    20 = @(8, 9)
    21 = @(7, 7)
    22 = @(7, 7)
    23 = @(8, 8)
    24 = @(7, 7)
    25 = @(7, 7)
    26 = @(7, 7)
    27 = @(9, 9)
    28 = @(8, 8)
    29 = @(7, 7)
    30 = @(7, 7)
    31 = @(7, 7)
    32 = @(8, 8)
    33 = @(7, 7)
    34 = @(9, 9)
    35 = @(9, 9)
    36 = @(7, 7)
    37 = @(8, 8)
    38 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
